$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the date in A1 (was 24-Apr-2024, now 24-May-2024)
$ws.Range("A1").Value = 45436

# Update prices
$ws.Range("D22").Value = 315
$ws.Range("D38").Value = 367.127

# Re-create the merged cell regions so their storage order matches
# the regenerated workbook (order-only change, same regions).
$ws.Range("B21:C21").UnMerge()
$ws.Range("B22:C22").UnMerge()
$ws.Range("A1:E1").UnMerge()
$ws.Range("B37:C37").UnMerge()
$ws.Range("B38:C38").UnMerge()

$ws.Range("B38:C38").Merge()
$ws.Range("B21:C21").Merge()
$ws.Range("B37:C37").Merge()
$ws.Range("B22:C22").Merge()
$ws.Range("A1:E1").Merge()
